$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list update -- values refreshed by scheduled scraper run

$ws.Range('D2').Value = '43.937.04'
$ws.Range('E2').Value = '  -0.85%  '
$ws.Range('D3').Value = '2.352.32'
$ws.Range('E3').Value = '  -0.80%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.674'
$ws.Range('E5').Value = '  -3.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '240.38'
$ws.Range('E6').Value = '  -1.35%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '73.32'
$ws.Range('E7').Value = '  -1.57%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.603'
$ws.Range('E9').Value = '  +0.56%  '
$ws.Range('E10').Value = '  -2.72%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '58.49'
$ws.Range('E11').Value = '  +1.29%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '33.72'
$ws.Range('E12').Value = '  +5.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.34'
$ws.Range('E13').Value = '  -2.44%  '
$ws.Range('E14').Value = '  -0.59%  '
$ws.Range('D15').Value = '2.704.61'
$ws.Range('E15').Value = '  -0.61%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '16.43'
$ws.Range('E16').Value = '  -4.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.912'
$ws.Range('E17').Value = '  -1.13%  '
$ws.Range('D18').Value = '2.351.19'
$ws.Range('E18').Value = '  -0.91%  '
$ws.Range('D19').Value = '43.830.59'
$ws.Range('E19').Value = '  -1.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000102'
$ws.Range('E20').Value = '  -1.75%  '
$ws.Range('E21').Value = '  -0.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '77.77'
$ws.Range('E22').Value = '  -1.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '257.05'
$ws.Range('E23').Value = '  -1.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.94'
$ws.Range('E24').Value = '  +15.64%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('B26').Value = 'WEMIXToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.75'
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('E27').Value = '  -2.93%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.65'
$ws.Range('E28').Value = '  -2.28%  '
$ws.Range('E29').Value = '  +2.71%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.67'
$ws.Range('E30').Value = '  -0.59%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '177.59'
$ws.Range('E31').Value = '  +1.64%  '
$ws.Range('E32').Value = '  -0.24%  '
$ws.Range('E33').Value = '  -0.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0757'
$ws.Range('E34').Value = '  -0.71%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.21'
$ws.Range('E35').Value = '  -3.87%  '
$ws.Range('E36').Value = '  +1.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.81'
$ws.Range('E37').Value = '  -2.83%  '
$ws.Range('B38').Value = 'THORChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.45'
$ws.Range('E38').Value = '  -2.11%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.39'
$ws.Range('E39').Value = '  -3.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0278'
$ws.Range('E40').Value = '  +0.15%  '
$ws.Range('B41').Value = 'MultiversX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '67.59'
$ws.Range('E41').Value = '  +26.00%  '
$ws.Range('B42').Value = 'FTXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.21'
$ws.Range('E42').Value = '  +16.85%  '
$ws.Range('E43').Value = '  +9.84%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '9.19'
$ws.Range('E44').Value = '  +0.46%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.05'
$ws.Range('E45').Value = '  -1.03%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.203'
$ws.Range('E46').Value = '  +2.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.51'
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('E48').Value = '  -0.82%  '
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('E50').Value = '  -3.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '99.25'
$ws.Range('E51').Value = '  -2.19%  '
